$wb = $excel.ActiveWorkbook

# --- Table 1 sheet ---
$t1 = $wb.Worksheets.Item("Table 1")

# Row 5: Age class -> Age class, years ; values become ranged bins
$t1.Range("A5").Value = "Age class, years"
$t1.Range("B5").Value = "16-30: 20% (n = 61)`n31-65: 66% (n = 202)`n>65: 14% (n = 44)"

# Row 15: Pre-existing somatic illness -> Pre-existing somatic illness type ; expand breakdown
$t1.Range("A15").Value = "Pre-existing somatic illness type"
$t1.Range("B15").Value = "none: 85% (n = 260)`nCVD: 2.9% (n = 9)`nneurological: 1.3% (n = 4)`nmetabolic: 1.3% (n = 4)`npulmonary: 0.65% (n = 2)`ncancer: 0.65% (n = 2)`nrheumatoid: 0.33% (n = 1)`nskin: 0.33% (n = 1)`nother: 7.8% (n = 24)"

# --- Table 2 sheet ---
$t2 = $wb.Worksheets.Item("Table 2")

# Row 3 (Sport type): mountain -> climbing/hiking/mountaineering
$t2.Range("B3").Value = "ski/snowboard: 64% (n = 197)`nsledding: 3.9% (n = 12)`nclimbing/hiking/mountaineering: 14% (n = 42)`nbiking: 16% (n = 48)`nother: 2.6% (n = 8)`nn = 307"

# --- Table 3 sheet ---
$t3 = $wb.Worksheets.Item("Table 3")

# Row 9 (EUROHIS-QOL 8 score): update statistic values
$t3.Range("B9").Value = "4.4 [IQR: 4 - 4.6]`nrange: 2 - 5"
